# The workbook tracks a "dry_fraction" column (D) that is being renamed to
# "wet_fraction". Updating the header cell's text is the whole semantic
# change; all the shared-string re-numbering seen in the diff falls out of
# that automatically once the old "dry_fraction" string is no longer
# referenced anywhere in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "wet_fraction"

# Match the author's final selection state (cell I14 was selected when the
# file was saved).
$ws.Range("I14").Select() | Out-Null
